$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (copy style/format from the existing last header cell AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every data row
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 76
    $ws.Cells.Item($row, 31).Value = 86
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Host "done"
